$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the female1/female2/female4/female5/female7/female8/female9/male1/male2
# rows (originally rows 2-10), shifting subsequent rows up.
$ws.Range("A2:B10").EntireRow.Delete() | Out-Null

# Reset the view: scroll back to the top and select A2:B28 with active cell A2.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A2:B28").Select() | Out-Null
